# 3513-RBI-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-PERIODIC-SP-FLAT-PENALTY-Newcreateloan.xlsx
# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gets a new (blank) column inserted before
# column N, pushing the old "Late" / "Outstanding" (heading) / "Outstanding"
# columns one slot to the right, widening the table from A1:P14 to A1:Q14.
# The sheet also becomes the active tab/selection of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before N (shifts N..P -> O..Q, clearing the new N).
$ws.Columns("N").Insert()

# Give the freshly inserted column the same width as column M (10.71 chars).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with K16 selected.
$ws.Activate()
$ws.Range("K16").Select()
